$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New store rows (Croma partner-brand stores) to append ----
$newRows = @(
    @("store_001018","Croma -Chhatrapati Sambhaji Nagar-Prozone Mall","Chhatrapati Sambhaji Nagar","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001019","Croma -Nashik-Solitario","Nashik","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001020","Croma -Baramati-Bhigwan Road","Baramati","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001021","Croma -Chhatrapati Sambhaji Nagar-Waluj","Chhatrapati Sambhaji Nagar","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001022","Croma -Nashik-Nashik Road","Nashik","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001023","Croma -Raipur-VIP Chowk","Raipur","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001024","Croma -Bilaspur-Srikant Verma Marg","Bilaspur","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001025","Croma -Solapur-Murarji Peth","Solapur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001026","Croma -Latur-Ambejogai Road","Latur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001027","Croma -Amravati-Badnera Road","AMRAVATI","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001028","Croma -Bhilai-Supela","Bhilai","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001029","Croma -Akola-Kirti Nagar","Akola","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001030","Croma -Ahmednagar-Kohinoor Mall","Ahmednagar","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001031","Croma -Nashik-Pathardi Phata","Nashik","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001032","Croma -Raipur-GE Road","Raipur","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001033","Croma -Chandrapur-Nagpur Road","Chandrapur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001034","Croma -Nagpur-Wardha Road","Nagpur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001035","Croma -Jalna-Old Mondha","Jalna","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001036","Croma -Rajnandgaon-Basantpur Road","Rajnandgaon","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001037","Croma -Ichalkranjii-Fortune Plaza II","Ichalkranjii","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001038","Croma -Solapur-Hotgi Road","Solapur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001039","Croma -Miraj-Vantmure Corner","Miraj","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001040","Croma -Chhatrapati Sambhaji Nagar-Jalna Road","Chhatrapati Sambhaji Nagar","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001041","Croma -Beed-Shivaji Maharaj Chowk","Beed","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001042","Croma -Nashik-Dindori Road","Nashik","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001043","Croma -Nashik-Gangapur Road","Nashik","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001044","Croma -Raipur-Bhatagaon","Raipur","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001045","Croma -Durg-Station Road","Durg","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001046","Croma -Korba-Transport Nagar","Korba","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001047","Croma -Nagpur-Manish Nagar","Nagpur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001048","Croma -Nanded-ITI Road","Nanded","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001049","Croma -Nagpur-Ganeshpeth","Nagpur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001050","Croma -Pune-Kothrud 2","Pune","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001051","Croma -Bilaspur-Sarkanada","Bilaspur","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001052","Croma -Ratnagiri-Arihant Mall","Ratnagiri","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001053","Croma -Sangamner-College Road","Sangamner","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001054","Croma -Nashik-Govind Nagar","Nashik","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001055","Croma -Machilipatnam-Ramanaidupet","Machilipatnam","brand_002","executive_00003","Vikash dubey"),
    @("store_001056","Croma -Erode-Perundurai Road","Erode","brand_002","executive_00003","Vikash dubey"),
    @("store_001057","Croma -Chennai-Tambaram West","Chennai","brand_002","executive_00003","Vikash dubey"),
    @("store_001058","Croma -Sangli-Ram Mandir Chowk","Sangli","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001059","Croma -Shahad-Kalyan Ahmednagar Highway","Shahad","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001060","Croma -Ranchi-Kanke Road","Ranchi","brand_002","executive_00016 ,executive_00003","Sanjay , Vikash Dubey"),
    @("store_001061","Croma -Kolhapur -Station Road","Kolhapur","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001062","Croma -Chhatrapati Sambhaji Nagar-Beed Bypass Road","Chhatrapati Sambhaji Nagar","brand_002","executive_00011 , executive_00018 ,executive_00003","Rushikesh , Soham , Vikash Dubey"),
    @("store_001063","Croma -Raipur-Pandri","Raipur","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey"),
    @("store_001064","Croma -Thane-Teen Hath Naka","Thane","brand_002","executive_00004 ,executive_00003","Soham , Vikash Dubey"),
    @("store_001065","Croma -Delhi-Airport Terminal 1 D","Delhi","brand_002","executive_00007 ,executive_00003","Ayush , Vikash Dubey"),
    @("store_001066","Croma -Guwahati ? Lokhara","Guwahati","brand_002","executive_00016 ,executive_00003","Sanjay , Vikash Dubey"),
    @("store_001067","Croma -Ghaziabad-Ambedkar Road","Ghaziabad","brand_002","executive_00002 ,executive_00003","Kanishk , Vikash Dubey")
)

$startRow = 1019

# Copy cell formatting from existing rows so new rows match the sheet's look:
#  - column A uses the highlighted "store id" style (copied from A718)
#  - columns B:F use the plain style (copied from row 2)
$ws.Range("A718").Copy() | Out-Null
$ws.Range("A$startRow`:A" + ($startRow + $newRows.Length - 1)).PasteSpecial(-4122) | Out-Null
$ws.Range("B2:F2").Copy() | Out-Null
$ws.Range("B$startRow`:F" + ($startRow + $newRows.Length - 1)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

$lastRow = $startRow + $newRows.Length - 1

# Re-fit the Store Name / City columns now that longer values were added
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Move the view / selection to the newly added data, like the author did
$ws.Range("B$lastRow").Select()
$excel.ActiveWindow.ScrollRow = 1055

# The defined name behind the (now stale) AutoFilter lost its target range
$names = $wb.Names
$names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!#REF!"
